# Update data parsing logic
# Appends a new data row (row 86) to each of the four worksheets,
# mirroring the existing "time"/"总长"/"ID"/"实际长度"/"和校验"/..._DEC
# table layout.

$wb = $excel.ActiveWorkbook

$rowData = @{
    "DE_LFT_#1" = @{
        A = 45872.43554398148
        B = "0x01,0x7c"
        C = "0x00,0xa6,0x46,0x93,0x3c,0x23,0x3f,0x43,0xe8,0xa0,"
        D = "0x01,0x34"
        E = "0x14"
        F = 380
        G = [double]"7.598631275147109e+23"
        H = 308
        I = 14
    }
    "DE_LFT_#2" = @{
        A = 45872.43554398148
        B = "0x01,0x7c"
        C = "0x00,0xa6,0x60,0x33,0x96,0x39,0x62,0xd0,0x5e,0x78,"
        D = "0x01,0x38"
        E = "0xe"
        F = 380
        G = [double]"5.68432987514711e+23"
        H = 312
        I = 14
    }
    "DE_PLT_#1" = @{
        A = 45872.43554398148
        B = "0x00,0x82"
        C = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x0b,0x40,0x0c,"
        D = "0x00,0x77"
        E = "0x7"
        F = 130
        G = [double]"5.68631262647114e+23"
        H = 119
        I = 7
    }
    "DE_PLT_#2" = @{
        A = 45872.43554398148
        B = "0x00,0x82"
        C = "0xd0,0x97,0x78,0x01,0x00,0x00,0x0e,0x3f,0x0c,0x0c,"
        D = "0x00,0x75"
        E = "0x3"
        F = 130
        G = [double]"9.85046333984776e+23"
        H = 117
        I = 3
    }
}

foreach ($ws in $wb.Worksheets) {
    $data = $rowData[$ws.Name]
    if ($data -eq $null) { continue }

    $r = $ws.UsedRange.Rows.Count + 1

    $ws.Cells.Item($r, 1).Value = $data.A
    $ws.Cells.Item($r, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"

    $ws.Cells.Item($r, 2).Value = $data.B
    $ws.Cells.Item($r, 3).Value = $data.C
    $ws.Cells.Item($r, 4).Value = $data.D
    $ws.Cells.Item($r, 5).Value = $data.E
    $ws.Cells.Item($r, 6).Value = $data.F
    $ws.Cells.Item($r, 7).Value = $data.G
    $ws.Cells.Item($r, 8).Value = $data.H
    $ws.Cells.Item($r, 9).Value = $data.I
}
